# Clear out the "Full Name" column data (A2:A13) while keeping the header
# (A1) and the cell formatting intact. This mirrors the author's change of
# emptying the sample "Full Name" values used to exercise the name-parsing
# macro, guarding against crashes when that column is blank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A13").ClearContents()

# Move the active selection, matching the post-edit workbook state.
$ws.Range("F7").Select()
